$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.85
$ws.Range("G2").Value = 5.4
$ws.Range("H2").Value = 1.79
$ws.Range("I2").Value = 2.02
$ws.Range("J2").Value = 3.25
$ws.Range("K2").Value = 4.6
$ws.Range("L2").Value = 1.39
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 3.15
$ws.Range("O2").Value = 1.31
$ws.Range("P2").Value = 1.76
$ws.Range("Q2").Value = 1.87
$ws.Range("R2").Value = 1.33
$ws.Range("S2").Value = 3.05
$ws.Range("T2").Value = 1.81
$ws.Range("U2").Value = 1.95
$ws.Range("V2").Value = 1.98
$ws.Range("W2").Value = 1.22
$ws.Range("X2").Value = 17
$ws.Range("Y2").Value = 10.5
$ws.Range("Z2").Value = 14
$ws.Range("AA2").Value = 25
$ws.Range("AB2").Value = 19.5
$ws.Range("AC2").Value = 10.5
$ws.Range("AD2").Value = 12.5
$ws.Range("AE2").Value = 23
$ws.Range("AF2").Value = 42
$ws.Range("AG2").Value = 22
$ws.Range("AH2").Value = 22
$ws.Range("AI2").Value = 42
$ws.Range("AJ2").Value = 120
$ws.Range("AK2").Value = 80
$ws.Range("AL2").Value = 70
$ws.Range("AM2").Value = 120
$ws.Range("AN2").Value = 95
$ws.Range("AO2").Value = 17

# Row 3
$ws.Range("N3").Value = 1.34
$ws.Range("O3").Value = 1.17
$ws.Range("P3").Value = 1.34
$ws.Range("Q3").Value = 1.17

# Row 4
$ws.Range("F4").Value = 3.95
$ws.Range("G4").Value = 5.9
$ws.Range("H4").Value = 1.71
$ws.Range("I4").Value = 2.06
$ws.Range("J4").Value = 3.1
$ws.Range("K4").Value = 4.6
$ws.Range("L4").Value = 1.4
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 3.05
$ws.Range("O4").Value = 1.32
$ws.Range("P4").Value = 1.8
$ws.Range("Q4").Value = 1.94
$ws.Range("R4").Value = 1.3
$ws.Range("S4").Value = 3.15
$ws.Range("T4").Value = 1.82
$ws.Range("U4").Value = 1.94
$ws.Range("V4").Value = 1.94
$ws.Range("W4").Value = 1.2
$ws.Range("X4").Value = 16.5
$ws.Range("Y4").Value = 10.5
$ws.Range("Z4").Value = 14
$ws.Range("AA4").Value = 27
$ws.Range("AB4").Value = 19
$ws.Range("AC4").Value = 10
$ws.Range("AD4").Value = 12.5
$ws.Range("AE4").Value = 26
$ws.Range("AF4").Value = 44
$ws.Range("AG4").Value = 23
$ws.Range("AH4").Value = 24
$ws.Range("AI4").Value = 48
$ws.Range("AJ4").Value = 140
$ws.Range("AK4").Value = 80
$ws.Range("AL4").Value = 85
$ws.Range("AN4").Value = 95
$ws.Range("AO4").Value = 17.5

# Row 6
$ws.Range("F6").Value = 2.12
$ws.Range("G6").Value = 2.62
$ws.Range("H6").Value = 2.84
$ws.Range("I6").Value = 3.9
$ws.Range("J6").Value = 2.92
$ws.Range("K6").Value = 4.3
$ws.Range("L6").Value = 1.38
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 3.05
$ws.Range("O6").Value = 1.3
$ws.Range("P6").Value = 1.73
$ws.Range("Q6").Value = 1.82
$ws.Range("R6").Value = 1.3
$ws.Range("S6").Value = 3.05
$ws.Range("T6").Value = 1.71
$ws.Range("U6").Value = 2
$ws.Range("V6").Value = 1.34
$ws.Range("W6").Value = 1.62
$ws.Range("X6").Value = 16.5
$ws.Range("Y6").Value = 15.5
$ws.Range("Z6").Value = 29
$ws.Range("AB6").Value = 12
$ws.Range("AC6").Value = 9.6
$ws.Range("AD6").Value = 18
$ws.Range("AE6").Value = 50
$ws.Range("AF6").Value = 18.5
$ws.Range("AG6").Value = 14
$ws.Range("AH6").Value = 22
$ws.Range("AI6").Value = 65
$ws.Range("AJ6").Value = 40
$ws.Range("AK6").Value = 32
$ws.Range("AL6").Value = 48
$ws.Range("AN6").Value = 25
$ws.Range("AO6").Value = 50

# Row 7
$ws.Range("L7").Value = 1.3
$ws.Range("R7").Value = 1.25

# Row 8
$ws.Range("F8").Value = 2.4
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 2.7
$ws.Range("I8").Value = 3.4
$ws.Range("J8").Value = 3.2
$ws.Range("L8").Value = 1.31
$ws.Range("M8").Value = 1.06
$ws.Range("N8").Value = 3.45
$ws.Range("O8").Value = 1.28
$ws.Range("P8").Value = 1.96
$ws.Range("Q8").Value = 1.73
$ws.Range("R8").Value = 1.38
$ws.Range("S8").Value = 2.78
$ws.Range("T8").Value = 1.66
$ws.Range("U8").Value = 2.16
$ws.Range("V8").Value = 1.43
$ws.Range("W8").Value = 1.56

# Row 9
$ws.Range("F9").Value = 1.6
$ws.Range("G9").Value = 1.61
$ws.Range("H9").Value = 5.5
$ws.Range("J9").Value = 4.8
$ws.Range("K9").Value = 5.2
$ws.Range("S9").Value = 2.22
$ws.Range("W9").Value = 2.62
$ws.Range("AA9").Value = 130
$ws.Range("AD9").Value = 22
$ws.Range("AE9").Value = 65
$ws.Range("AJ9").Value = 1000
$ws.Range("AO9").Value = 980

# Row 10
$ws.Range("G10").Value = 2.96
$ws.Range("H10").Value = 2.48
$ws.Range("I10").Value = 2.7
$ws.Range("Q10").Value = 1.74
$ws.Range("S10").Value = 2.66
$ws.Range("W10").Value = 1.51

# Row 11
$ws.Range("F11").Value = 3.4
$ws.Range("J11").Value = 3.8
$ws.Range("L11").Value = 1.28
$ws.Range("N11").Value = 5.5
$ws.Range("Q11").Value = 1.58
$ws.Range("R11").Value = 1.62
$ws.Range("S11").Value = 2.36
$ws.Range("T11").Value = 1.54
$ws.Range("U11").Value = 2.6
$ws.Range("Y11").Value = 15
$ws.Range("AC11").Value = 9.800000000000001
$ws.Range("AO11").Value = 11.5

# Row 12
$ws.Range("F12").Value = 6
$ws.Range("G12").Value = 6.6
$ws.Range("H12").Value = 1.6
$ws.Range("I12").Value = 1.63
$ws.Range("J12").Value = 4.3
$ws.Range("K12").Value = 4.7
$ws.Range("P12").Value = 2.44
$ws.Range("Q12").Value = 1.65
$ws.Range("R12").Value = 1.58
$ws.Range("S12").Value = 2.5
$ws.Range("T12").Value = 1.71
$ws.Range("U12").Value = 2.26
$ws.Range("V12").Value = 2.58
$ws.Range("W12").Value = 1.17
$ws.Range("X12").Value = 22
$ws.Range("Y12").Value = 11.5
$ws.Range("Z12").Value = 11
$ws.Range("AA12").Value = 16
$ws.Range("AB12").Value = 28
$ws.Range("AC12").Value = 10.5
$ws.Range("AE12").Value = 15.5
$ws.Range("AF12").Value = 55
$ws.Range("AG12").Value = 23
$ws.Range("AH12").Value = 18.5
$ws.Range("AI12").Value = 980
$ws.Range("AJ12").Value = 210
$ws.Range("AK12").Value = 75
$ws.Range("AL12").Value = 70
$ws.Range("AM12").Value = 85
$ws.Range("AN12").Value = 85
$ws.Range("AO12").Value = 6.6
